$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9223046214701632
$ws.Range("J2").Value = 0.9223046214701632
$ws.Range("O2").Value = 0.8416031693647025
$ws.Range("P2").Value = 0.8416031693647025
$ws.Range("S2").Value = 0.7762144925490015
$ws.Range("T2").Value = 0.7762144925490015

# Row 3
$ws.Range("I3").Value = 0.9223046214701632
$ws.Range("J3").Value = 0.9223046214701632
$ws.Range("M3").Value = 0.2972526666666667
$ws.Range("N3").Value = 0.8917580000000001
$ws.Range("O3").Value = 0.1583968306352975
$ws.Range("P3").Value = 0.1583968306352975
$ws.Range("Q3").Value = 0.4959723166393334
$ws.Range("R3").Value = 4.463750849754001
$ws.Range("S3").Value = 0.1460901289211616
$ws.Range("T3").Value = 0.1460901289211616

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.140557
$ws.Range("H4").Value = 0.421671
$ws.Range("I4").Value = 0.07769537852983674
$ws.Range("J4").Value = 0.07769537852983674
$ws.Range("O4").Value = 0.8416031693647025
$ws.Range("P4").Value = 0.8416031693647025
$ws.Range("Q4").Value = 0.22199291466
$ws.Range("R4").Value = 1.99793623194
$ws.Range("S4").Value = 0.06538867681570086
$ws.Range("T4").Value = 0.06538867681570086

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.140557
$ws.Range("H5").Value = 0.421671
$ws.Range("I5").Value = 0.07769537852983674
$ws.Range("J5").Value = 0.07769537852983674
$ws.Range("M5").Value = 0.2972526666666667
$ws.Range("N5").Value = 0.8917580000000001
$ws.Range("O5").Value = 0.1583968306352975
$ws.Range("P5").Value = 0.1583968306352975
$ws.Range("Q5").Value = 0.04178094306866667
$ws.Range("R5").Value = 0.3760284876180001
$ws.Range("S5").Value = 0.01230670171413588
$ws.Range("T5").Value = 0.01230670171413588
